$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price (D) and Volume(1h) (E) columns so the
# numeric-looking / percent-looking strings are stored as literal text,
# matching the inline-string cells already in the workbook.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "311.87"
$ws.Range("E2").Value = "1.76%"
$ws.Range("D3").Value = "37.67"
$ws.Range("E3").Value = "0.51%"
$ws.Range("D4").Value = "5.128"
$ws.Range("E4").Value = "0.88%"
$ws.Range("D5").Value = "0.07880"
$ws.Range("E5").Value = "1.96%"
$ws.Range("D6").Value = "4.412"
$ws.Range("E6").Value = "1.74%"
$ws.Range("D7").Value = "1.906"
$ws.Range("E7").Value = "1.15%"
$ws.Range("E8").Value = "1.14%"
$ws.Range("D9").Value = "2.857"
$ws.Range("E9").Value = "-6.87%"
$ws.Range("D10").Value = "0.9182"
$ws.Range("E10").Value = "-0.58%"
$ws.Range("D11").Value = "0.1179"
$ws.Range("E11").Value = "-3.99%"
$ws.Range("D12").Value = "0.1930"
$ws.Range("E12").Value = "3.07%"
$ws.Range("D13").Value = "0.09082"
$ws.Range("E13").Value = "3.63%"
$ws.Range("D14").Value = "0.03319"
$ws.Range("E14").Value = "-2.51%"
$ws.Range("D15").Value = "0.09600"
$ws.Range("E15").Value = "-1.09%"
$ws.Range("D16").Value = "0.001383"
$ws.Range("E16").Value = "0.94%"
$ws.Range("D17").Value = "0.006017"
$ws.Range("E17").Value = "-0.98%"
$ws.Range("D18").Value = "3.541"
$ws.Range("E18").Value = "-1.28%"
$ws.Range("E19").Value = "0.96%"
$ws.Range("D20").Value = "5.284"
$ws.Range("E20").Value = "5.33%"
$ws.Range("D21").Value = "0.1284"
$ws.Range("E21").Value = "1.25%"
$ws.Range("E22").Value = "3.93%"
$ws.Range("D23").Value = "0.04355"
$ws.Range("E23").Value = "0.60%"
$ws.Range("E24").Value = "3.14%"
$ws.Range("E25").Value = "10.40%"
$ws.Range("E26").Value = "0.66%"
$ws.Range("D27").Value = "0.0003986"
$ws.Range("E27").Value = "-98.11%"
$ws.Range("D39").Value = "0.02263"
$ws.Range("E39").Value = "3.78%"
$ws.Range("D40").Value = "0.05097"
$ws.Range("E40").Value = "4.13%"
$ws.Range("D41").Value = "0.007442"
$ws.Range("E41").Value = "-0.55%"
$ws.Range("D42").Value = "0.009036"
$ws.Range("E42").Value = "-8.66%"
$ws.Range("E43").Value = "1.12%"
$ws.Range("D44").Value = "0.001998"
$ws.Range("E44").Value = "-2.99%"
$ws.Range("D45").Value = "0.008597"
$ws.Range("E45").Value = "-12.60%"
$ws.Range("D46").Value = "0.00006562"
$ws.Range("E46").Value = "0.39%"
$ws.Range("E47").Value = "-0.05%"
$ws.Range("D48").Value = "0.0009992"
$ws.Range("E48").Value = "-23.18%"
$ws.Range("D49").Value = "0.003017"
$ws.Range("E49").Value = "0.72%"
$ws.Range("E50").Value = "-0.05%"
$ws.Range("E51").Value = "-0.05%"
